# Update PLC data 2025-10-13 14:00:32
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 171543
$ws.Range("C4").Value = 162345
$ws.Range("C7").Value = 5.36
$ws.Range("C8").Value = 65.83
